$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.716399999999998
$ws.Range("B6").Value = 5.351700000000005
$ws.Range("B7").Value = 5.182299999999998
$ws.Range("C7").Value = -14.17459999999999
$ws.Range("C12").Value = -10.86639999999999
$ws.Range("E12").Value = 17.54790000000001
$ws.Range("D13").Value = -8.534900000000002
$ws.Range("D14").Value = -8.126100000000001
$ws.Range("C15").Value = -14.75389999999999
$ws.Range("B16").Value = 7.556999999999995
$ws.Range("D16").Value = -8.571000000000005
$ws.Range("D19").Value = -8.777099999999992
$ws.Range("B20").Value = 9.367999999999997
$ws.Range("C20").Value = -12.0203
$ws.Range("C21").Value = -11.95920000000002
$ws.Range("C22").Value = -12.4082
$ws.Range("D22").Value = -8.124900000000006
$ws.Range("E22").Value = 16.56100000000001
$ws.Range("C23").Value = -12.17760000000001
$ws.Range("B28").Value = 6.024000000000002
$ws.Range("B29").Value = 5.346699999999999
$ws.Range("C29").Value = -11.59110000000001
$ws.Range("E29").Value = 17.04640000000001
$ws.Range("B32").Value = 7.705599999999993
$ws.Range("C34").Value = -11.80770000000001
$ws.Range("E34").Value = 17.3973
$ws.Range("D36").Value = -8.463999999999997
$ws.Range("B40").Value = 9.565699999999994
$ws.Range("C42").Value = -12.54470000000001
$ws.Range("C43").Value = -12.53269999999999
$ws.Range("E43").Value = 17.23290000000001
$ws.Range("C44").Value = -13.2354
$ws.Range("C45").Value = -13.7245
$ws.Range("B46").Value = 6.180499999999996
$ws.Range("C46").Value = -13.1947
$ws.Range("D46").Value = -8.060999999999996
$ws.Range("E48").Value = 17.47820000000001
$ws.Range("C50").Value = -13.82399999999999
$ws.Range("D50").Value = -8.064599999999999
$ws.Range("B51").Value = 6.257400000000002
$ws.Range("C51").Value = -11.96700000000001
$ws.Range("B52").Value = 5.372099999999998
$ws.Range("B57").Value = 5.438599999999996
$ws.Range("B59").Value = 5.539800000000001
$ws.Range("E60").Value = 16.5977
$ws.Range("B62").Value = 5.4392
$ws.Range("B66").Value = 5.910499999999997
$ws.Range("C66").Value = -11.86510000000001
$ws.Range("C67").Value = -11.78779999999999
$ws.Range("E68").Value = 17.55150000000001
$ws.Range("E70").Value = 18.38060000000002
$ws.Range("B73").Value = 8.058400000000002
$ws.Range("E73").Value = 17.6312
$ws.Range("B74").Value = 9.244499999999993
$ws.Range("C79").Value = -12.4552
$ws.Range("C84").Value = -13.34099999999999
$ws.Range("E87").Value = 16.33339999999999
$ws.Range("B92").Value = 5.571599999999997
$ws.Range("C92").Value = -11.44429999999999
$ws.Range("E92").Value = 18.34950000000001
$ws.Range("D95").Value = -8.037699999999999
$ws.Range("C97").Value = -12.04620000000001
$ws.Range("D97").Value = -8.671099999999996
$ws.Range("B100").Value = 5.740799999999997
$ws.Range("E101").Value = 16.81330000000001
